$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows (2-296)
$ws.Range("C2:C296").Value = 45184

# 2. Row 296 gets an explicit row height (15pt, custom height) in the new version
$ws.Rows.Item(296).RowHeight = 15

# 3. Add the new row 297 with the same formatting as row 296, then set its values
$ws.Range("A296:R296").Copy()
$ws.Range("A297:R297").PasteSpecial(-4122)

$ws.Cells.Item(297, 1).Value = "A 42527-2023"
$ws.Cells.Item(297, 2).Value = 45181
$ws.Cells.Item(297, 3).Value = 45184
$ws.Cells.Item(297, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(297, 5).Value = "VANSBRO"
$ws.Cells.Item(297, 6).Value = ""
$ws.Cells.Item(297, 7).Value = 14.3
$ws.Cells.Item(297, 8).Value = 0
$ws.Cells.Item(297, 9).Value = 0
$ws.Cells.Item(297, 10).Value = 0
$ws.Cells.Item(297, 11).Value = 0
$ws.Cells.Item(297, 12).Value = 0
$ws.Cells.Item(297, 13).Value = 0
$ws.Cells.Item(297, 14).Value = 0
$ws.Cells.Item(297, 15).Value = 0
$ws.Cells.Item(297, 16).Value = 0
$ws.Cells.Item(297, 17).Value = 0
